$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs / Wnt2 / Fzd3 / ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.265744
$ws.Range("N2").Value = 0.7972319999999999
$ws.Range("O2").Value = 0.1034864391735229
$ws.Range("P2").Value = 0.1034864391735229
$ws.Range("Q2").Value = 0.669889601152
$ws.Range("R2").Value = 6.029006410368
$ws.Range("S2").Value = 0.1034864391735229
$ws.Range("T2").Value = 0.1034864391735229

# Row 3 (FAPs / Wnt2 / Fzd3 / FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.307583
$ws.Range("N3").Value = 0.922749
$ws.Range("O3").Value = 0.1197794472135201
$ws.Range("P3").Value = 0.1197794472135201
$ws.Range("Q3").Value = 0.775357687064
$ws.Range("R3").Value = 6.978219183576
$ws.Range("S3").Value = 0.1197794472135201
$ws.Range("T3").Value = 0.1197794472135201

# Row 4 (FAPs / Wnt2 / Fzd3 / sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.994584333333333
$ws.Range("N4").Value = 5.983753
$ws.Range("O4").Value = 0.7767341136129571
$ws.Range("P4").Value = 0.7767341136129571
$ws.Range("Q4").Value = 5.027964144141333
$ws.Range("R4").Value = 45.251677297272
$ws.Range("S4").Value = 0.7767341136129571
$ws.Range("T4").Value = 0.7767341136129571
